$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$values_ALC = @{
    "H17" = 536.59186
    "J17" = 543.6042
    "L17" = 1630.8126
    "N17" = -1966.8126
    "H137" = 3751352.2
    "I137" = 1725244.2
    "J137" = 9092909
    "K137" = 5175732.6
    "L137" = 27278727
    "M137" = -5173182.6
    "N137" = -27283827
}
foreach ($key in $values_ALC.Keys) {
    $ws.Range($key).Value = $values_ALC[$key]
}

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$values_ARM = @{
    "H32" = 32217.959
    "I32" = 6883.657
    "K32" = 6883.657
    "M32" = -6596.657
    "H44" = 37666
    "J44" = 37666
    "L44" = 37666
    "N44" = -38642
    "H55" = 23323.572
    "J55" = 23323.572
    "L55" = 23323.572
    "N55" = -23953.572
    "H61" = 1982.84
    "I61" = 1823.55
    "J61" = 2620
    "K61" = 1823.55
    "L61" = 2620
    "M61" = -1611.55
    "N61" = -3044
    "H80" = 31554.666
    "J80" = 31554.666
    "L80" = 31554.666
    "N80" = -33550.666
    "H83" = 31554.666
    "J83" = 31554.666
    "L83" = 94663.99800000001
    "N83" = -104647.998
    "H132" = 149769.44
    "I132" = 180076.75
    "J132" = 8335.333000000001
    "K132" = 540230.25
    "L132" = 25005.999
    "M132" = -537700.25
    "N132" = -30065.999
    "H136" = 1982.84
    "I136" = 1823.55
    "J136" = 2620
    "K136" = 5470.65
    "L136" = 7860
    "M136" = -2920.65
    "N136" = -12960
}
foreach ($key in $values_ARM.Keys) {
    $ws.Range($key).Value = $values_ARM[$key]
}

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$values_BSM = @{
    "H82" = 8419
    "I82" = 8419
    "K82" = 8419
    "M82" = -8036
    "H85" = 8419
    "I85" = 8419
    "K85" = 8419
    "M85" = -7093
}
foreach ($key in $values_BSM.Keys) {
    $ws.Range($key).Value = $values_BSM[$key]
}

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$values_CRP = @{
    "H31" = 1809.0984
    "I31" = 1621.05
    "J31" = 2167.2856
    "K31" = 1621.05
    "L31" = 2167.2856
    "M31" = -1326.05
    "N31" = -2757.2856
    "H34" = 1809.0984
    "I34" = 1621.05
    "J34" = 2167.2856
    "K34" = 1621.05
    "L34" = 2167.2856
    "M34" = -1419.05
    "N34" = -2571.2856
    "H41" = 15741.25
    "I41" = 5000
    "J41" = 19321.666
    "K41" = 5000
    "L41" = 19321.666
    "M41" = -4572
    "N41" = -20177.666
    "H50" = 6702.125
    "J50" = 6945.2856
    "L50" = 6945.2856
    "N50" = -8195.285599999999
    "H60" = 11103
    "I60" = 0
    "J60" = 11103
    "K60" = 0
    "L60" = 11103
    "N60" = -12125
    "H109" = 20514
    "J109" = 20514
    "L109" = 20514
    "N109" = -22594
    "H132" = 2254.4119
    "I132" = 1673.52
    "J132" = 3868
    "K132" = 5020.559999999999
    "L132" = 11604
    "M132" = -2490.559999999999
    "N132" = -16664
}
foreach ($key in $values_CRP.Keys) {
    $ws.Range($key).Value = $values_CRP[$key]
}
foreach ($ref in @("M60")) {
    $ws.Range($ref).ClearContents()
}

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$values_CUL = @{
    "H5" = 351456.25
    "I5" = 394.5
    "J5" = 1112090.1
    "K5" = 1183.5
    "L5" = 3336270.3
    "M5" = -1071.5
    "N5" = -3336494.3
    "H107" = 713.6957
    "I107" = 337.1
    "J107" = 1003.38464
    "K107" = 1011.3
    "L107" = 3010.15392
    "M107" = 908.6999999999999
    "N107" = -6850.15392
    "H122" = 46645.477
    "I122" = 286.44446
    "J122" = 51672.363
    "K122" = 2578.00014
    "L122" = 465051.267
    "M122" = -128.0001400000001
    "N122" = -469951.267
    "H135" = 351456.25
    "I135" = 394.5
    "J135" = 1112090.1
    "K135" = 3550.5
    "L135" = 10008810.9
    "M135" = -1015.5
    "N135" = -10013880.9
}
foreach ($key in $values_CUL.Keys) {
    $ws.Range($key).Value = $values_CUL[$key]
}

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$values_GSM = @{
    "H57" = 13928.25
    "J57" = 13928.25
    "L57" = 13928.25
    "N57" = -15568.25
    "H70" = 41942.855
    "I70" = 49756.523
    "J70" = 6000
    "K70" = 49756.523
    "L70" = 6000
    "M70" = -49486.523
    "N70" = -6540
    "H73" = 41942.855
    "I73" = 49756.523
    "J73" = 6000
    "K73" = 49756.523
    "L73" = 6000
    "M73" = -48820.523
    "N73" = -7872
    "H132" = 2702.1667
    "I132" = 2424.3333
    "J132" = 3813.5
    "K132" = 7272.999899999999
    "L132" = 11440.5
    "M132" = -4742.999899999999
    "N132" = -16500.5
}
foreach ($key in $values_GSM.Keys) {
    $ws.Range($key).Value = $values_GSM[$key]
}

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$values_LTW = @{
    "H14" = 1640
    "J14" = 2750
    "L14" = 2750
    "N14" = -3094
    "H21" = 11625
    "I21" = 4500
    "J21" = 14000
    "K21" = 4500
    "L21" = 14000
    "M21" = -4326
    "N21" = -14348
    "H55" = 768.1053000000001
    "I55" = 449.57144
    "J55" = 953.9167
    "K55" = 449.57144
    "L55" = 953.9167
    "M55" = -276.57144
    "N55" = -1299.9167
    "H109" = 28483.334
    "J109" = 28483.334
    "L109" = 28483.334
    "N109" = -31257.334
    "H132" = 2122.7222
    "I132" = 1445.1111
    "J132" = 2800.3333
    "K132" = 4335.3333
    "L132" = 8400.999899999999
    "M132" = -1805.3333
    "N132" = -13460.9999
    "H133" = 41057.184
    "J133" = 41057.184
    "L133" = 41057.184
    "N133" = -46117.184
    "H136" = 1700.4884
    "I136" = 1326.4839
    "J136" = 2666.6667
    "K136" = 3979.4517
    "L136" = 8000.000100000001
    "M136" = -1429.4517
    "N136" = -13100.0001
}
foreach ($key in $values_LTW.Keys) {
    $ws.Range($key).Value = $values_LTW[$key]
}

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$values_WVR = @{
    "H25" = 15999.667
    "J25" = 15999.667
    "L25" = 15999.667
    "N25" = -16585.667
    "H109" = 0
    "J109" = 0
    "L109" = 0
    "H136" = 1317.9111
    "I136" = 1138.15
    "K136" = 3414.45
    "M136" = -864.4500000000003
}
foreach ($key in $values_WVR.Keys) {
    $ws.Range($key).Value = $values_WVR[$key]
}
foreach ($ref in @("N109")) {
    $ws.Range($ref).ClearContents()
}
